# msz - field hint and error checks part 1
# Adds a new row (23) to Tabelle1 describing a check for hints regarding
# mandatory fields on the Vehicle page, mirroring the existing
# "Vehicle Page check for open mandatory fields" / "Insurant Page check
# for open mandatory fields" rows already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$hintText = "Vehicle Page check for hints regarding mandatory fields"

$ws.Range("A23").Value = $hintText
$ws.Range("B23").Value = "<CHK>"
$ws.Range("C23").Value = $hintText
$ws.Range("H23").Value = "<NOP>"

# Move the active selection to the newly added row, matching the saved
# workbook's cursor position.
$ws.Range("A23").Select() | Out-Null
